# This edit re-orders the data rows (2..21) of the sheet: each destination
# row ends up with the Fecha/Variedad/Calidad/Volumen/Precio.../Unidad/
# Origen/Precio $/Kg values (columns D, K, L, M, N, O, P, Q, R, S) that used
# to live in a different source row, while columns A, B, C, E, F, G, H, I, J
# and T (which are constant across all rows) are left untouched.
#
# Mapping: destination row -> source row (values are read from the ORIGINAL
# workbook state before any writes happen).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 16
    3  = 12
    4  = 4
    5  = 17
    6  = 10
    7  = 5
    8  = 6
    9  = 14
    10 = 8
    11 = 7
    12 = 13
    13 = 2
    14 = 19
    15 = 18
    16 = 20
    17 = 3
    18 = 9
    19 = 21
    20 = 11
    21 = 15
}

$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S")

# 1. Snapshot the original values for the columns that move, for every row
#    that participates in the shuffle (2..21), before any writes occur.
$snapshot = @{}
foreach ($r in 2..21) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2. Write the snapshotted values into their destination rows according to
#    the mapping.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
